$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "78÷6="  # was 96÷6=
$t.Cell(1, 2).Range.Text = "66÷7="  # was 90÷6=
$t.Cell(1, 3).Range.Text = "88÷7="  # was 30÷8=
$t.Cell(1, 4).Range.Text = "23÷6="  # was 90÷6=
$t.Cell(1, 5).Range.Text = "24÷2="  # was 30÷9=

$t.Cell(5, 1).Range.Text = "13÷6="  # was 97÷8=
$t.Cell(5, 2).Range.Text = "32÷5="  # was 62÷4=
$t.Cell(5, 3).Range.Text = "64÷8="  # was 44÷3=
$t.Cell(5, 4).Range.Text = "79÷8="  # was 53÷5=
$t.Cell(5, 5).Range.Text = "84÷8="  # was 11÷6=

$t.Cell(9, 1).Range.Text = "82÷3="  # was 67÷4=
$t.Cell(9, 2).Range.Text = "72÷7="  # was 10÷2=
$t.Cell(9, 3).Range.Text = "17÷9="  # was 37÷7=
$t.Cell(9, 4).Range.Text = "73÷6="  # was 93÷7=
$t.Cell(9, 5).Range.Text = "27÷3="  # was 17÷4=

$t.Cell(13, 1).Range.Text = "95÷5="  # was 90÷9=
$t.Cell(13, 2).Range.Text = "24÷2="  # was 30÷8=
$t.Cell(13, 3).Range.Text = "19÷2="  # was 19÷4=
$t.Cell(13, 4).Range.Text = "58÷6="  # was 41÷4=
$t.Cell(13, 5).Range.Text = "25÷6="  # was 20÷3=

$t.Cell(17, 1).Range.Text = "51÷5="  # was 20÷2=
$t.Cell(17, 2).Range.Text = "86÷7="  # was 47÷8=
$t.Cell(17, 3).Range.Text = "26÷2="  # was 37÷2=
$t.Cell(17, 4).Range.Text = "79÷9="  # was 15÷4=
$t.Cell(17, 5).Range.Text = "74÷6="  # was 74÷7=
